$wb = $excel.ActiveWorkbook

# Rename sheets (renaming updates workbook.xml <sheet name="..."> entries)
$wb.Worksheets.Item(1).Name = "GNG_TO-16504777930765538"
$wb.Worksheets.Item(2).Name = "NB_TO-16504777953505538"
$wb.Worksheets.Item(3).Name = "RS_TO-16504777953515549"
$wb.Worksheets.Item(4).Name = "TOL_TO-1650477795398555"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16504777954595876"

# Sheet1 (GNG_TO) - update B2:B5
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16504777930395548.csv"
$ws1.Range("B3").Value = "GNG_stims-16504777930595922.csv"
$ws1.Range("B4").Value = "go_stims-16504777930615566.csv"
$ws1.Range("B5").Value = "GNG_stims-16504777930755835.csv"

# Sheet2 (NB_TO) - update B2:B10
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "ZB-match_5-1650477793178556.csv"
$ws2.Range("B3").Value = "TB-16504777950195868.csv"
$ws2.Range("B4").Value = "OB-1650477793469557.csv"
$ws2.Range("B5").Value = "ZB-match_4-1650477793207558.csv"
$ws2.Range("B6").Value = "TB-16504777950565875.csv"
$ws2.Range("B7").Value = "TB-16504777953395584.csv"
$ws2.Range("B8").Value = "OB-16504777936855876.csv"
$ws2.Range("B9").Value = "ZB-match_8-16504777932425532.csv"
$ws2.Range("B10").Value = "OB-16504777935865536.csv"

# Sheet4 (TOL_TO) - update B2:B7
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16504777953665552.csv"
$ws4.Range("B3").Value = "ZM_stims-16504777953535547.csv"
$ws4.Range("B4").Value = "MM_stims-16504777953825555.csv"
$ws4.Range("B5").Value = "ZM_stims-16504777953675532.csv"
$ws4.Range("B6").Value = "MM_stims-16504777953975809.csv"
$ws4.Range("B7").Value = "ZM_stims-16504777953825555.csv"

# Sheet5 (vSAT_TO) - update B2:B5
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16504777954015532.csv"
$ws5.Range("B3").Value = "vSAT_stims-16504777954285896.csv"
$ws5.Range("B4").Value = "SAT_stims-16504777954125607.csv"
$ws5.Range("B5").Value = "vSAT_stims-16504777954455543.csv"
